$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place character edits to preserve rich-text runs) ---

# A8: "Volume 30   Number  9" -> "Volume 30   Number  10"  (the trailing "9" -> "10")
$volCell = $ws.Range("A8")
$volChars = $volCell.Characters(21, 1)
$volChars.Text = "10"

# C9: "Report Covering the Week  2/27/2023  Through  3/5/2023"
#  -> "Report Covering the Week  3/6/2023  Through  3/12/2023"
# Replace the second (rightmost) date first so the first date's character
# offset is unaffected by the length change of the replacement text.
$dateCell = $ws.Range("C9")
$dateChars2 = $dateCell.Characters(47, 8)
$dateChars2.Text = "3/12/2023"
$dateChars1 = $ws.Range("C9").Characters(27, 9)
$dateChars1.Text = "3/6/2023"

# --- Weekly crime statistics table updates (rows 14-30, cols C-N) ---

# Row 14
$ws.Range("C14").Value = 6
$ws.Range("E14").Value = -25
$ws.Range("F14").Value = 21
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = -36.363636363636
$ws.Range("I14").Value = 66
$ws.Range("J14").Value = 81
$ws.Range("K14").Value = -18.518518518518
$ws.Range("L14").Value = -17.5
$ws.Range("M14").Value = -22.352941176470
$ws.Range("N14").Value = -82.539682539682

# Row 15
$ws.Range("C15").Value = 33
$ws.Range("D15").Value = 37
$ws.Range("E15").Value = -10.810810810810
$ws.Range("F15").Value = 111
$ws.Range("G15").Value = 150
$ws.Range("H15").Value = -26
$ws.Range("I15").Value = 291
$ws.Range("J15").Value = 334
$ws.Range("K15").Value = -12.874251497006
$ws.Range("L15").Value = 18.292682926829
$ws.Range("M15").Value = 32.272727272727
$ws.Range("N15").Value = -50.927487352445

# Row 16
$ws.Range("D16").Value = 284
$ws.Range("E16").Value = -5.281690140845
$ws.Range("F16").Value = 1162
$ws.Range("G16").Value = 1284
$ws.Range("H16").Value = -9.501557632398
$ws.Range("I16").Value = 2967
$ws.Range("J16").Value = 3036
$ws.Range("K16").Value = -2.272727272727
$ws.Range("L16").Value = 41.690544412607
$ws.Range("M16").Value = -13.042203985932
$ws.Range("N16").Value = -82.214362786236

# Row 17
$ws.Range("C17").Value = 454
$ws.Range("D17").Value = 433
$ws.Range("E17").Value = 4.849884526558
$ws.Range("F17").Value = 1832
$ws.Range("G17").Value = 1718
$ws.Range("H17").Value = 6.635622817229
$ws.Range("I17").Value = 4716
$ws.Range("J17").Value = 4270
$ws.Range("K17").Value = 10.444964871194
$ws.Range("L17").Value = 32.583637897104
$ws.Range("M17").Value = 61.617546264564
$ws.Range("N17").Value = -30.288248337028

# Row 18
$ws.Range("C18").Value = 270
$ws.Range("D18").Value = 322
$ws.Range("E18").Value = -16.149068322981
$ws.Range("F18").Value = 1074
$ws.Range("G18").Value = 1229
$ws.Range("H18").Value = -12.611879576891
$ws.Range("I18").Value = 2820
$ws.Range("J18").Value = 2960
$ws.Range("K18").Value = -4.729729729729
$ws.Range("L18").Value = 21.394748170469
$ws.Range("M18").Value = -19.749573136027
$ws.Range("N18").Value = -85.507246376811

# Row 19
$ws.Range("C19").Value = 873
$ws.Range("E19").Value = -5.723542116630
$ws.Range("F19").Value = 3567
$ws.Range("G19").Value = 3746
$ws.Range("H19").Value = -4.778430325680
$ws.Range("I19").Value = 9129
$ws.Range("J19").Value = 9495
$ws.Range("K19").Value = -3.854660347551
$ws.Range("L19").Value = 57.641167328613
$ws.Range("M19").Value = 38.675375968403
$ws.Range("N19").Value = -38.706861823553

# Row 20
$ws.Range("C20").Value = 265
$ws.Range("D20").Value = 239
$ws.Range("E20").Value = 10.878661087866
$ws.Range("F20").Value = 1086
$ws.Range("G20").Value = 1013
$ws.Range("H20").Value = 7.206317867719
$ws.Range("I20").Value = 2789
$ws.Range("J20").Value = 2644
$ws.Range("K20").Value = 5.484114977307
$ws.Range("L20").Value = 101.663051337672
$ws.Range("M20").Value = 55.810055865921
$ws.Range("N20").Value = -87.513431232091

# Row 21
$ws.Range("C21").Value = 2170
$ws.Range("D21").Value = 2249
$ws.Range("E21").Value = -3.512672298799
$ws.Range("F21").Value = 8853
$ws.Range("G21").Value = 9173
$ws.Range("H21").Value = -3.488498855336
$ws.Range("I21").Value = 22778
$ws.Range("J21").Value = 22820
$ws.Range("K21").Value = -0.184049079754
$ws.Range("L21").Value = 47.201757787256
$ws.Range("M21").Value = 22.978080120937
$ws.Range("N21").Value = -71.915764555026

# Row 22
$ws.Range("C22").Value = 48
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = 14.285714285714
$ws.Range("F22").Value = 176
$ws.Range("G22").Value = 195
$ws.Range("H22").Value = -9.743589743589
$ws.Range("I22").Value = 391
$ws.Range("J22").Value = 470
$ws.Range("K22").Value = -16.808510638297
$ws.Range("L22").Value = 53.937007874015
$ws.Range("M22").Value = -1.511335012594

# Row 23
$ws.Range("C23").Value = 80
$ws.Range("D23").Value = 105
$ws.Range("E23").Value = -23.809523809523
$ws.Range("F23").Value = 426
$ws.Range("G23").Value = 422
$ws.Range("H23").Value = 0.947867298578
$ws.Range("I23").Value = 1117
$ws.Range("J23").Value = 1070
$ws.Range("K23").Value = 4.392523364485
$ws.Range("L23").Value = 22.343921139101
$ws.Range("M23").Value = 57.545839210155

# Row 24
$ws.Range("C24").Value = 2097
$ws.Range("D24").Value = 2083
$ws.Range("E24").Value = 0.672107537205
$ws.Range("F24").Value = 8249
$ws.Range("G24").Value = 8423
$ws.Range("H24").Value = -2.065772290157
$ws.Range("I24").Value = 20593
$ws.Range("J24").Value = 19652
$ws.Range("K24").Value = 4.788316710767
$ws.Range("L24").Value = 40.912823320104
$ws.Range("M24").Value = 48.193724812895

# Row 25
$ws.Range("C25").Value = 771
$ws.Range("D25").Value = 781
$ws.Range("E25").Value = -1.280409731113
$ws.Range("F25").Value = 3055
$ws.Range("G25").Value = 2995
$ws.Range("H25").Value = 2.003338898163
$ws.Range("I25").Value = 7666
$ws.Range("J25").Value = 7185
$ws.Range("K25").Value = 6.694502435629
$ws.Range("L25").Value = 39.154111453984
$ws.Range("M25").Value = -4.115071919949

# Row 26
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 65
$ws.Range("E26").Value = -32.307692307692
$ws.Range("G26").Value = 247
$ws.Range("H26").Value = -30.769230769230
$ws.Range("I26").Value = 460
$ws.Range("J26").Value = 531
$ws.Range("K26").Value = -13.370998116760
$ws.Range("L26").Value = 11.650485436893

# Row 27
$ws.Range("C27").Value = 111
$ws.Range("D27").Value = 99
$ws.Range("E27").Value = 12.121212121212
$ws.Range("F27").Value = 377
$ws.Range("G27").Value = 358
$ws.Range("H27").Value = 5.307262569832
$ws.Range("I27").Value = 934
$ws.Range("J27").Value = 843
$ws.Range("K27").Value = 10.794780545670
$ws.Range("L27").Value = 29.722222222222

# Row 28
$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 26
$ws.Range("E28").Value = -7.692307692307
$ws.Range("F28").Value = 70
$ws.Range("G28").Value = 77
$ws.Range("H28").Value = -9.090909090909
$ws.Range("I28").Value = 203
$ws.Range("J28").Value = 234
$ws.Range("K28").Value = -13.247863247863
$ws.Range("L28").Value = -2.870813397129
$ws.Range("M28").Value = -13.617021276595
$ws.Range("N28").Value = -80.518234165067

# Row 29
$ws.Range("C29").Value = 18
$ws.Range("D29").Value = 23
$ws.Range("E29").Value = -21.739130434782
$ws.Range("F29").Value = 57
$ws.Range("G29").Value = 71
$ws.Range("H29").Value = -19.718309859154
$ws.Range("I29").Value = 166
$ws.Range("J29").Value = 210
$ws.Range("K29").Value = -20.952380952381
$ws.Range("L29").Value = -13.089005235602
$ws.Range("M29").Value = -15.736040609137
$ws.Range("N29").Value = -82.617801047120

# Row 30
$ws.Range("D30").Value = 19
$ws.Range("E30").Value = -89.473684210526
$ws.Range("F30").Value = 22
$ws.Range("G30").Value = 92
$ws.Range("H30").Value = -76.086956521739
$ws.Range("I30").Value = 66
$ws.Range("J30").Value = 161
$ws.Range("K30").Value = -59.006211180124
$ws.Range("L30").Value = 29.411764705882
